$d = $word.ActiveDocument

# Step 1: remove the old _GoBack bookmark. In the original file it sat
# alone inside its own empty paragraph right after the title paragraph
# ("Testleitungsmanual"); that paragraph becomes a plain empty paragraph.
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

# Step 2: fill in the placeholder "#" with the actual figure "18" in the
# "Ihre Meinung ... Minuten gedauert." bullet point.
$d.Content.Find.Execute("ca. # Minuten", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ca. 18 Minuten", 2)

# Step 3: locate the paragraph that now holds the updated sentence.
$target = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*abgefragt*") {
        $target = $i
    }
}
$p = $d.Paragraphs.Item($target)
$text = $p.Range.Text
$paraStart = $p.Range.Start

# Step 4a: break the single run into two right after "...Fragebogen hat".
# Bookmarks are siblings of <w:r> in the OOXML, so adding (and promptly
# removing) one at this spot forces the surrounding text to split into
# separate runs, exactly like the diff shows.
$hatEnd = $paraStart + $text.IndexOf("hat") + ("hat").Length
$d.Bookmarks.Add("TempSplit", $d.Range($hatEnd, $hatEnd))
$d.Bookmarks.Item("TempSplit").Delete()

# Step 4b: re-create the _GoBack bookmark right after "ca. 18" (i.e. right
# before " Minuten gedauert."), matching where the edit actually occurred.
$ca18End = $paraStart + $text.IndexOf("ca. 18") + ("ca. 18").Length
$d.Bookmarks.Add("_GoBack", $d.Range($ca18End, $ca18End))
